$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Test Cases" sheet: add two new test script rows (PUBLONS013 / PUBLONS014)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Test Cases")

# Row 12 - PUBLONS013 (cells filled left to right; each new distinct string
# becomes a new shared-string table entry in this order)
$ws.Range("A12").Value = 'PUBLONS013'
$ws.Range("B12").Value = 'OPQA-5771||OPQA-5772||OPQA-5773||OPQA-5774||OPQA-5775||OPQA-5776||OPQA-5777||OPQA-6000'
$ws.Range("C12").Value = 'Verify Password must have at least one special character from !@#$%^*()~`{}[] in Registration  page||Verify Password must contain at least one number is ALWAYS enforced in Registration  page||Verify Password must have at least one alphabet character either upper or lower case is ALWAYS enforced in Registration  page||Verify that the Password minimum length of 8 characters is ALWAYS enforced in Registration  page.||Verify Password Maximum Length of 95 characters is ALWAYS enforced in Registration  page||Verify that error message "Password is too long" whenever enter more than 95 characters||Verify that "View password rules on the right" error message at the time of entering password||Verify that "Should not have leading and trailing spaces" error message at the time of entering password'
$ws.Range("D12").Value = 'Y'
$ws.Rows.Item(12).RowHeight = 90

# Row 13 - PUBLONS014 (B13/C13 filled before A13, matching the shared-string
# append order captured in the target workbook)
$ws.Range("B13").Value = 'OPQA-5770'
$ws.Range("C13").Value = 'Verify that error message "Please enter a password." whenever not enter any text in email field'
$ws.Range("D13").Value = 'Y'
$ws.Range("A13").Value = 'PUBLONS014'

# Column B got wider to fit the new content
$ws.Columns.Item(2).ColumnWidth = 74.67

# Make "Test Cases" the active tab / active sheet, with D13 selected.
# (Selecting a range on a sheet activates that sheet, so do this LAST -
#  "PUBLONS005", which used to be the active tab, keeps its own A4
#  selection untouched and simply stops being the active tab.)
$ws.Activate()
$ws.Range("D13").Select()
